# feat: add 2022-Q1 data
#
# 1) Insert a new sheet "2022-Q1" (holding the per-fund breakdown for the
#    quarter) right before the "总计" (totals) summary sheet.
# 2) Insert a new top data row into "总计" for "2022-Q1" (holding count=4,
#    value=0.91), pushing the existing quarters down by one row and
#    renumbering the helper index column (A).
#
# NOTE: sheet references returned by Worksheets.Item(...) can go stale once
# the Worksheets collection is structurally changed (Add/Copy/Move/Delete),
# so we always re-fetch a sheet by name right after such an operation
# instead of reusing an earlier variable.

$wb = $excel.ActiveWorkbook

# --- 1) Build the new "2022-Q1" sheet -------------------------------------
# Copy "2021-Q4" (same 8-column fund-breakdown layout/styles we need) so the
# header row / index-column formatting (style "s=2") comes along for free,
# then drop it in right before "总计".
$totalBeforeCopy = $wb.Worksheets.Item("总计")
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($totalBeforeCopy)

$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# 2021-Q4 has 10 data rows (rows 2-11); 2022-Q1 only needs 4 (rows 2-5), so
# drop the extra copied rows.
$newSheet.Range("A6:H11").Delete()

# A never-touched cell on the new sheet: used purely as a "no explicit
# style" formatting source further down so text-forcing the fund
# codes/numbers below doesn't leave a stray style index behind.
$blank = $newSheet.Range("Z100")

# Force text storage on the code / numeric-looking text columns so values
# such as "007066" or "3.30" keep their leading/trailing zeros instead of
# being auto-coerced to numbers.
$newSheet.Range("B2:B5").NumberFormat = "@"
$newSheet.Range("D2:G5").NumberFormat = "@"

$fundRows = @(
    @("213001", "宝盈鸿利收益灵活配置混合A", "17.98", "90.37", "3.84", "0.6904", 8),
    @("007066", "浦银安盛先进制造混合A",     "3.66",  "74.55", "3.30", "0.1208", 10),
    @("007067", "浦银安盛先进制造混合C",     "2.27",  "74.55", "3.30", "0.0749", 10),
    @("007581", "宝盈鸿利收益灵活配置混合C", "0.73",  "90.37", "3.84", "0.0280", 8)
)

$r = 2
foreach ($row in $fundRows) {
    $newSheet.Range("B$r").Value = $row[0]
    $newSheet.Range("C$r").Value = $row[1]
    $newSheet.Range("D$r").Value = $row[2]
    $newSheet.Range("E$r").Value = $row[3]
    $newSheet.Range("F$r").Value = $row[4]
    $newSheet.Range("G$r").Value = $row[5]
    $newSheet.Range("H$r").Value = $row[6]
    $r++
}

# Reapply the "no explicit style" formatting (copied from the untouched
# $blank cell) on top of the text values we just wrote, so the cells end up
# identical to their siblings (plain text cells, no style index) instead of
# keeping the temporary "@" number-format style.
$blank.Copy()
$newSheet.Range("B2:B5").PasteSpecial(-4122)
$newSheet.Range("D2:G5").PasteSpecial(-4122)

# --- 2) Insert the "2022-Q1" summary row into "总计" -----------------------
# Re-fetch by name: the Worksheets collection changed above (Copy), so any
# earlier reference to "总计" may now point at the wrong sheet.
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

# Restore the index-column style ("s=2") on the freshly inserted A2 by
# copying it from the row right below (still carrying the original style).
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.91

# Renumber the helper index column (A) for every row beneath the new one.
for ($row = 3; $row -le 7; $row++) {
    $total.Range("A$row").Value = $row - 2
}

# Restore the originally-active first tab (copying a sheet makes the copy
# the active tab).
$wb.Worksheets.Item(1).Activate()

Write-Output "2022-Q1 sheet added; zongji (totals) sheet updated"
